$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the greeting text in cell E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active/selected cell, matching the saved selection state
$ws.Activate()
$ws.Range("E8").Select()
